$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5646
$ws.Range("J40").Value = 7271.2856
$ws.Range("L40").Value = 7271.2856
$ws.Range("N40").Value = -7621.2856
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H87").Value = 54000
$ws.Range("J87").Value = 54000
$ws.Range("L87").Value = 54000
$ws.Range("N87").Value = -56496
$ws.Range("H88").Value = 4065.8572
$ws.Range("I88").Value = 3003
$ws.Range("K88").Value = 3003
$ws.Range("M88").Value = -2597
$ws.Range("H90").Value = 54000
$ws.Range("J90").Value = 54000
$ws.Range("L90").Value = 162000
$ws.Range("N90").Value = -174480
$ws.Range("H91").Value = 4065.8572
$ws.Range("I91").Value = 3003
$ws.Range("K91").Value = 3003
$ws.Range("M91").Value = -1599
$ws.Range("H112").Value = 3911.818
$ws.Range("J112").Value = 3903
$ws.Range("L112").Value = 11709
$ws.Range("N112").Value = -13925
$ws.Range("H138").Value = 2042.0426
$ws.Range("I138").Value = 1341.7368
$ws.Range("K138").Value = 4025.2104
$ws.Range("M138").Value = 1114.7896

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 24704
$ws.Range("J55").Value = 32032
$ws.Range("L55").Value = 32032
$ws.Range("N55").Value = -32662
$ws.Range("H74").Value = 2705.1155
$ws.Range("I74").Value = 1963.6052
$ws.Range("K74").Value = 1963.6052
$ws.Range("M74").Value = -1089.6052
$ws.Range("H77").Value = 2705.1155
$ws.Range("I77").Value = 1963.6052
$ws.Range("K77").Value = 9818.026
$ws.Range("M77").Value = -5450.026
$ws.Range("H88").Value = 2292
$ws.Range("I88").Value = 2500
$ws.Range("J88").Value = 2257.3333
$ws.Range("K88").Value = 2500
$ws.Range("L88").Value = 2257.3333
$ws.Range("M88").Value = -2094
$ws.Range("N88").Value = -3069.3333
$ws.Range("H91").Value = 2292
$ws.Range("I91").Value = 2500
$ws.Range("J91").Value = 2257.3333
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 2257.3333
$ws.Range("M91").Value = -1096
$ws.Range("N91").Value = -5065.3333
$ws.Range("H124").Value = 38500
$ws.Range("J124").Value = 38500
$ws.Range("L124").Value = 38500
$ws.Range("N124").Value = -48320

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 61066.332
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H20").Value = 3670.1304
$ws.Range("I20").Value = 3115.1
$ws.Range("K20").Value = 3115.1
$ws.Range("M20").Value = -2868.1
$ws.Range("H35").Value = 42434.832
$ws.Range("J35").Value = 45761.8
$ws.Range("L35").Value = 45761.8
$ws.Range("N35").Value = -46381.8
$ws.Range("H82").Value = 41400
$ws.Range("J82").Value = 41400
$ws.Range("L82").Value = 41400
$ws.Range("N82").Value = -42166
$ws.Range("H85").Value = 41400
$ws.Range("J85").Value = 41400
$ws.Range("L85").Value = 41400
$ws.Range("N85").Value = -44052
$ws.Range("H140").Value = 73499.5
$ws.Range("J140").Value = 73499.5
$ws.Range("L140").Value = 73499.5
$ws.Range("N140").Value = -83859.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 23099.5
$ws.Range("J41").Value = 22000
$ws.Range("L41").Value = 22000
$ws.Range("N41").Value = -22856
$ws.Range("H51").Value = 21000
$ws.Range("J51").Value = 27500
$ws.Range("L51").Value = 27500
$ws.Range("N51").Value = -28972
$ws.Range("H59").Value = 40000
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H61").Value = 21000
$ws.Range("J61").Value = 27500
$ws.Range("L61").Value = 27500
$ws.Range("N61").Value = -28196
$ws.Range("H68").Value = 41677.25
$ws.Range("J68").Value = 41677.25
$ws.Range("L68").Value = 41677.25
$ws.Range("N68").Value = -43175.25
$ws.Range("H71").Value = 41677.25
$ws.Range("J71").Value = 41677.25
$ws.Range("L71").Value = 125031.75
$ws.Range("N71").Value = -132519.75
$ws.Range("H74").Value = 41573.715
$ws.Range("J74").Value = 41573.715
$ws.Range("L74").Value = 41573.715
$ws.Range("N74").Value = -43321.715
$ws.Range("H77").Value = 41573.715
$ws.Range("J77").Value = 41573.715
$ws.Range("L77").Value = 124721.145
$ws.Range("N77").Value = -133457.145
$ws.Range("H99").Value = 3600
$ws.Range("I99").Value = 3777.6
$ws.Range("J99").Value = 3244.8
$ws.Range("K99").Value = 3777.6
$ws.Range("L99").Value = 3244.8
$ws.Range("M99").Value = -2279.6
$ws.Range("N99").Value = -6240.8
$ws.Range("H126").Value = 3600
$ws.Range("I126").Value = 3777.6
$ws.Range("J126").Value = 3244.8
$ws.Range("K126").Value = 11332.8
$ws.Range("L126").Value = 9734.400000000001
$ws.Range("M126").Value = -8862.799999999999
$ws.Range("N126").Value = -14674.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5999.5
$ws.Range("I56").Value = 5999.5
$ws.Range("K56").Value = 5999.5
$ws.Range("M56").Value = -5469.5
$ws.Range("H139").Value = 2089.842
$ws.Range("I139").Value = 1860.5333
$ws.Range("J139").Value = 2949.75
$ws.Range("K139").Value = 5581.5999
$ws.Range("L139").Value = 8849.25
$ws.Range("M139").Value = -441.5999000000002
$ws.Range("N139").Value = -19129.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 45644.168
$ws.Range("J46").Value = 51621.668
$ws.Range("L46").Value = 51621.668
$ws.Range("N46").Value = -51933.668
$ws.Range("H97").Value = 911.8
$ws.Range("I97").Value = 715.375
$ws.Range("J97").Value = 1697.5
$ws.Range("K97").Value = 715.375
$ws.Range("L97").Value = 1697.5
$ws.Range("M97").Value = -219.375
$ws.Range("N97").Value = -2689.5
$ws.Range("H99").Value = 18620
$ws.Range("I99").Value = 12160
$ws.Range("K99").Value = 12160
$ws.Range("M99").Value = -9914
$ws.Range("H132").Value = 2750
$ws.Range("I132").Value = 2666.6667
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 8000.000100000001
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -5470.000100000001
$ws.Range("N132").Value = -14060
$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H101").Value = 64489.5
$ws.Range("J101").Value = 64489.5
$ws.Range("L101").Value = 64489.5
$ws.Range("N101").Value = -70979.5
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5275.067
$ws.Range("J62").Value = 6445
$ws.Range("L62").Value = 6445
$ws.Range("N62").Value = -7693
$ws.Range("H65").Value = 5275.067
$ws.Range("J65").Value = 6445
$ws.Range("L65").Value = 32225
$ws.Range("N65").Value = -38465
$ws.Range("H81").Value = 777.7143
$ws.Range("I81").Value = 777.7143
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1555.4286
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -494.4286
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 777.7143
$ws.Range("I84").Value = 777.7143
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 7777.143
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -2473.143
$ws.Range("N84").ClearContents()
$ws.Range("H103").Value = 47000
$ws.Range("J103").Value = 47000
$ws.Range("L103").Value = 47000
$ws.Range("N103").Value = -49344
